$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2576.2632
$ws.Range("J17").Value = 2996.8667
$ws.Range("L17").Value = 8990.6001
$ws.Range("N17").Value = -9326.6001

$ws.Range("H64").Value = 4371.2856
$ws.Range("J64").Value = 3149.75
$ws.Range("L64").Value = 3149.75
$ws.Range("N64").Value = -3645.75

$ws.Range("H67").Value = 4371.2856
$ws.Range("J67").Value = 3149.75
$ws.Range("L67").Value = 3149.75
$ws.Range("N67").Value = -4865.75

$ws.Range("H92").Value = 1324.6666
$ws.Range("I92").Value = 1123.7778
$ws.Range("J92").Value = 1626
$ws.Range("K92").Value = 1123.7778
$ws.Range("L92").Value = 1626
$ws.Range("M92").Value = 124.2221999999999
$ws.Range("N92").Value = -4122

$ws.Range("H98").Value = 914.1
$ws.Range("I98").Value = 815.6667
$ws.Range("J98").Value = 1800
$ws.Range("K98").Value = 815.6667
$ws.Range("L98").Value = 1800
$ws.Range("M98").Value = 682.3333
$ws.Range("N98").Value = -4796

$ws.Range("H101").Value = 370
$ws.Range("I101").Value = 370
$ws.Range("K101").Value = 1110
$ws.Range("M101").Value = 512

$ws.Range("H122").Value = 914.1
$ws.Range("I122").Value = 815.6667
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 2447.0001
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = 2.999899999999798
$ws.Range("N122").Value = -10300

$ws.Range("H132").Value = 1242.875
$ws.Range("I132").Value = 1176.963
$ws.Range("K132").Value = 3530.889
$ws.Range("M132").Value = -1000.889

$ws.Range("H138").Value = 4525.254
$ws.Range("I138").Value = 1995.5714
$ws.Range("K138").Value = 5986.7142
$ws.Range("M138").Value = -846.7142000000003

$ws.Range("H141").Value = 6038
$ws.Range("I141").Value = 5797.5
$ws.Range("K141").Value = 17392.5
$ws.Range("M141").Value = -12212.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 1877.7778
$ws.Range("J33").Value = 2500
$ws.Range("L33").Value = 2500
$ws.Range("N33").Value = -3158

$ws.Range("H45").Value = 3373.3333
$ws.Range("I45").Value = 1482.75
$ws.Range("K45").Value = 1482.75
$ws.Range("M45").Value = -1105.75

$ws.Range("H51").Value = 24000
$ws.Range("J51").Value = 24000
$ws.Range("L51").Value = 24000
$ws.Range("N51").Value = -25512

$ws.Range("H61").Value = 2994.8572
$ws.Range("I61").Value = 2993
$ws.Range("K61").Value = 2993
$ws.Range("M61").Value = -2781

$ws.Range("H132").Value = 1811.8
$ws.Range("I132").Value = 1727.2142
$ws.Range("K132").Value = 5181.642599999999
$ws.Range("M132").Value = -2651.642599999999

$ws.Range("H136").Value = 2994.8572
$ws.Range("I136").Value = 2993
$ws.Range("K136").Value = 8979
$ws.Range("M136").Value = -6429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 634.4167
$ws.Range("I22").Value = 634.4167
$ws.Range("K22").Value = 634.4167
$ws.Range("M22").Value = -461.4167

$ws.Range("H86").Value = 3076.25
$ws.Range("I86").Value = 499.5
$ws.Range("J86").Value = 5653
$ws.Range("K86").Value = 499.5
$ws.Range("L86").Value = 5653
$ws.Range("M86").Value = 623.5
$ws.Range("N86").Value = -7899

$ws.Range("H89").Value = 3076.25
$ws.Range("I89").Value = 499.5
$ws.Range("J89").Value = 5653
$ws.Range("K89").Value = 2497.5
$ws.Range("L89").Value = 28265
$ws.Range("M89").Value = 3118.5
$ws.Range("N89").Value = -39497

$ws.Range("H94").Value = 707.875
$ws.Range("I94").Value = 707.5714
$ws.Range("J94").Value = 710
$ws.Range("K94").Value = 707.5714
$ws.Range("L94").Value = 710
$ws.Range("M94").Value = -256.5714
$ws.Range("N94").Value = -1612

$ws.Range("H97").Value = 10789.667
$ws.Range("I97").Value = 10789.667
$ws.Range("K97").Value = 10789.667
$ws.Range("M97").Value = -9798.666999999999

$ws.Range("H105").Value = 3297.75
$ws.Range("I105").Value = 2670.6296
$ws.Range("J105").Value = 4600.231
$ws.Range("K105").Value = 2670.6296
$ws.Range("L105").Value = 4600.231
$ws.Range("M105").Value = -923.6296000000002
$ws.Range("N105").Value = -8094.231

$ws.Range("H122").Value = 424999.2
$ws.Range("I122").Value = 299999
$ws.Range("J122").Value = 508332.66
$ws.Range("K122").Value = 299999
$ws.Range("L122").Value = 508332.66
$ws.Range("M122").Value = -295099
$ws.Range("N122").Value = -518132.66

$ws.Range("H134").Value = 4948.3
$ws.Range("I134").Value = 4898.1113
$ws.Range("K134").Value = 14694.3339
$ws.Range("M134").Value = -12159.3339

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 8975.333000000001
$ws.Range("I122").Value = 8956
$ws.Range("K122").Value = 26868
$ws.Range("M122").Value = -24418

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 91514950
$ws.Range("I4").Value = 105963600
$ws.Range("J4").Value = 6851.3335
$ws.Range("K4").Value = 317890800
$ws.Range("L4").Value = 20554.0005
$ws.Range("M4").Value = -317890688
$ws.Range("N4").Value = -20778.0005

$ws.Range("H60").Value = 733
$ws.Range("I60").Value = 618.25
$ws.Range("J60").Value = 962.5
$ws.Range("K60").Value = 1854.75
$ws.Range("L60").Value = 2887.5
$ws.Range("M60").Value = -1603.75
$ws.Range("N60").Value = -3389.5

$ws.Range("H68").Value = 1701.8889
$ws.Range("I68").Value = 2129.3333
$ws.Range("J68").Value = 1488.1666
$ws.Range("K68").Value = 6387.999899999999
$ws.Range("L68").Value = 4464.4998
$ws.Range("M68").Value = -5576.999899999999
$ws.Range("N68").Value = -6086.4998

$ws.Range("H71").Value = 1701.8889
$ws.Range("I71").Value = 2129.3333
$ws.Range("J71").Value = 1488.1666
$ws.Range("K71").Value = 19163.9997
$ws.Range("L71").Value = 13393.4994
$ws.Range("M71").Value = -15107.9997
$ws.Range("N71").Value = -21505.4994

$ws.Range("H121").Value = 1420

$ws.Range("H129").Value = 2772.4614
$ws.Range("I129").Value = 2007.8334
$ws.Range("J129").Value = 3427.8572
$ws.Range("K129").Value = 6023.5002
$ws.Range("L129").Value = 10283.5716
$ws.Range("M129").Value = -1023.5002
$ws.Range("N129").Value = -20283.5716

$ws.Range("H131").Value = 3227.35
$ws.Range("I131").Value = 3006
$ws.Range("K131").Value = 9018
$ws.Range("M131").Value = -3978

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4156.75
$ws.Range("J132").Value = 4923.6665
$ws.Range("L132").Value = 14770.9995
$ws.Range("N132").Value = -19830.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2577.25
$ws.Range("I7").Value = 2577.25
$ws.Range("K7").Value = 2577.25
$ws.Range("M7").Value = -2465.25

$ws.Range("H16").Value = 5568.2856
$ws.Range("I16").Value = 3885.3333
$ws.Range("K16").Value = 3885.3333
$ws.Range("M16").Value = -3715.3333

$ws.Range("H126").Value = 2577.25
$ws.Range("I126").Value = 2577.25
$ws.Range("K126").Value = 7731.75
$ws.Range("M126").Value = -5261.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1205.1613
$ws.Range("I113").Value = 912.7368
$ws.Range("K113").Value = 2738.2104
$ws.Range("M113").Value = -568.2103999999999

$ws.Range("H122").Value = 1583.9166
$ws.Range("I122").Value = 1583.9166
$ws.Range("K122").Value = 4751.7498
$ws.Range("M122").Value = -2301.7498

$ws.Range("H126").Value = 79686.234
$ws.Range("I126").Value = 112770.78
$ws.Range("K126").Value = 338312.34
$ws.Range("M126").Value = -335842.34
